$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 10 (B10/C10): long PT objectives text -> short "Robson" text
$ws.Range("B10:C10").Value = '7455355 - Robson da Silva Rocha'

# Remove old rows 13-24 entirely (content + row-specific heights), then rebuild 13-23
$ws.Range("A13:C24").EntireRow.Delete()

# Row 13
$ws.Rows.Item(13).Insert()
$ws.Range("A13").Value = 'Programa resumido:'
$ws.Range("B13").Value = 'Semestral'
$ws.Range("C13").Value = 'Semestral'
$ws.Rows.Item(13).RowHeight = 60

# Row 14
$ws.Rows.Item(14).Insert()
$ws.Range("A14").Value = 'Short syllabus:'
$ws.Range("B14").Value = 'Introduction to qualitative analysis, indicating its applicability and limitations. Use of qualitative techniques to analyze the main ions of environmental importance. Analysis of solids, particles, sediments. Studies of samples of environmental importance.'
$ws.Range("C14").Value = 'Introduction to qualitative analysis, indicating its applicability and limitations. Use of qualitative techniques to analyze the main ions of environmental importance. Analysis of solids, particles, sediments. Studies of samples of environmental importance.'
$ws.Rows.Item(14).RowHeight = 60

# Row 15
$ws.Rows.Item(15).Insert()
$ws.Range("A15").Value = 'Programa:'
$ws.Range("B15").Value = '01/01/2022'
$ws.Range("C15").Value = '01/01/2022'
$ws.Rows.Item(15).RowHeight = 120

# Row 16
$ws.Rows.Item(16).Insert()
$ws.Range("A16").Value = 'Syllabus:'
$ws.Range("B16").Value = '- Review of laboratory safety rules- Introduction to qualitative analysis: Definitions, objectives and limitations.- Analysis of solids, particles, sediments.- Identification of group I cations (K+, Na+ and NH4+); group II (Mg2+, Ca2+ and Ba2+); group III (Al3+, Fe3+, Fe2+, Mn2+).- Study of anions and their applications in environmental analysis (Cl and its species, SO42-, CO32-, S2-, NO3-).- Gravimetric analysis: fundamentals and precipitate formation.- Analysis of the main cations and anions in known and unknown samples for students- Analysis of metals in soil, water or other important environmental samples'
$ws.Range("C16").Value = '- Review of laboratory safety rules- Introduction to qualitative analysis: Definitions, objectives and limitations.- Analysis of solids, particles, sediments.- Identification of group I cations (K+, Na+ and NH4+); group II (Mg2+, Ca2+ and Ba2+); group III (Al3+, Fe3+, Fe2+, Mn2+).- Study of anions and their applications in environmental analysis (Cl and its species, SO42-, CO32-, S2-, NO3-).- Gravimetric analysis: fundamentals and precipitate formation.- Analysis of the main cations and anions in known and unknown samples for students- Analysis of metals in soil, water or other important environmental samples'
$ws.Rows.Item(16).RowHeight = 120

# Row 17
$ws.Rows.Item(17).Insert()
$ws.Range("A17").Value = 'Avaliação:'

# Row 18
$ws.Rows.Item(18).Insert()
$ws.Range("A18").Value = 'Método:'
$ws.Range("B18").Value = '7455355 - Robson da Silva Rocha'
$ws.Range("C18").Value = '7455355 - Robson da Silva Rocha'
$ws.Rows.Item(18).RowHeight = 60

# Row 19
$ws.Rows.Item(19).Insert()
$ws.Range("A19").Value = 'Critério:'
$ws.Range("B19").Value = 'O método de avaliação será composto por avaliações teóricas, práticas e relatórios de atividades de práticas laboratoriais.'
$ws.Range("C19").Value = 'O método de avaliação será composto por avaliações teóricas, práticas e relatórios de atividades de práticas laboratoriais.'
$ws.Rows.Item(19).RowHeight = 60

# Row 20
$ws.Rows.Item(20).Insert()
$ws.Range("A20").Value = 'Norma de recuperação:'
$ws.Range("B20").Value = 'Para o cálculo da nota final (NF) será feita a média aritmética das avaliações aplicadas. Estará aprovado por notas o aluno que obtiver nota final igual ou superior a 5,0 pontos.'
$ws.Range("C20").Value = 'Para o cálculo da nota final (NF) será feita a média aritmética das avaliações aplicadas. Estará aprovado por notas o aluno que obtiver nota final igual ou superior a 5,0 pontos.'
$ws.Rows.Item(20).RowHeight = 60

# Row 21
$ws.Rows.Item(21).Insert()
$ws.Range("A21").Value = 'Bibliografia:'
$ws.Range("B21").Value = 'Avaliação de recuperação (R) envolvendo todo o conteúdo da disciplina. Média Final = (NF+R) / 2 => 5,0 Aprovado'
$ws.Range("C21").Value = 'Avaliação de recuperação (R) envolvendo todo o conteúdo da disciplina. Média Final = (NF+R) / 2 => 5,0 Aprovado'
$ws.Rows.Item(21).RowHeight = 120

# Row 22
$ws.Rows.Item(22).Insert()
$ws.Range("A22").Value = 'Requisitos:'

# Row 23
$ws.Rows.Item(23).Insert()
$ws.Range("B23").Value = 'LOQ4098 -  Fundamentos de Química para Engenharia II (Requisito fraco)
'
$ws.Range("C23").Value = 'LOQ4098 -  Fundamentos de Química para Engenharia II (Requisito fraco)
'
$ws.Rows.Item(23).RowHeight = 30
